# Update "paises.xlsx" with the latest COVID-19 snapshot (countries & provincias Spain)
# - refreshes the "Datos actualizados a ..." timestamp banner
# - refreshes numeric columns for several countries whose counts moved
# - Nigeria overtook Armenia in total cases, and Lesoto overtook Gibraltar/
#   Bermudas/Camboya/Brunei, so those rows swap rank (and therefore swap the
#   country name shown on that row) while keeping the table sorted by
#   "Casos totales" (column B) descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp banner (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 11 de Julio de 2020 a las 01:04"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 3285506
$ws.Cells.Item(4, 3).Value = 65507
$ws.Cells.Item(4, 4).Value = 1453158
$ws.Cells.Item(4, 5).Value = 1695778
$ws.Cells.Item(4, 7).Value = 748
$ws.Cells.Item(4, 8).Value = 136570

# Row 22: Colombia
$ws.Cells.Item(22, 2).Value = 140776
$ws.Cells.Item(22, 3).Value = 6803
$ws.Cells.Item(22, 4).Value = 58800
$ws.Cells.Item(22, 5).Value = 77051
$ws.Cells.Item(22, 7).Value = 211
$ws.Cells.Item(22, 8).Value = 4925

# Row 23: Canada
$ws.Cells.Item(23, 2).Value = 107125
$ws.Cells.Item(23, 3).Value = 320
$ws.Cells.Item(23, 4).Value = 70901
$ws.Cells.Item(23, 5).Value = 27465

# Row 25: Argentina
$ws.Cells.Item(25, 2).Value = 94060
$ws.Cells.Item(25, 3).Value = 3367
$ws.Cells.Item(25, 5).Value = 53302
$ws.Cells.Item(25, 7).Value = 54
$ws.Cells.Item(25, 8).Value = 1774

# Row 34: Kazajistan
$ws.Cells.Item(34, 4).Value = 31815
$ws.Cells.Item(34, 5).Value = 22668

# Row 52: was Armenia, now Nigeria (Nigeria's total overtakes Armenia's)
$ws.Cells.Item(52, 1).Value = "Nigeria"
$ws.Cells.Item(52, 2).Value = 31323
$ws.Cells.Item(52, 3).Value = 575
$ws.Cells.Item(52, 4).Value = 12795
$ws.Cells.Item(52, 5).Value = 17819
$ws.Cells.Item(52, 7).Value = 20
$ws.Cells.Item(52, 8).Value = 709

# Row 53: was Nigeria, now Armenia (keeps Armenia's previous totals, one row lower)
$ws.Cells.Item(53, 1).Value = "Armenia"
$ws.Cells.Item(53, 2).Value = 30903
$ws.Cells.Item(53, 3).Value = 557
$ws.Cells.Item(53, 4).Value = 18709
$ws.Cells.Item(53, 5).Value = 11648
$ws.Cells.Item(53, 7).Value = 11
$ws.Cells.Item(53, 8).Value = 546

# Row 54: Guatemala
$ws.Cells.Item(54, 2).Value = 27619
$ws.Cells.Item(54, 3).Value = 961
$ws.Cells.Item(54, 4).Value = 4024
$ws.Cells.Item(54, 5).Value = 22456
$ws.Cells.Item(54, 7).Value = 47
$ws.Cells.Item(54, 8).Value = 1139

# Row 59: Japon
$ws.Cells.Item(59, 2).Value = 20719
$ws.Cells.Item(59, 3).Value = 348
$ws.Cells.Item(59, 4).Value = 17652
$ws.Cells.Item(59, 5).Value = 2085
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 982

# Row 142: Uruguay
$ws.Cells.Item(142, 2).Value = 985
$ws.Cells.Item(142, 3).Value = 8
$ws.Cells.Item(142, 4).Value = 886

# Row 150: Surinam
$ws.Cells.Item(150, 4).Value = 481
$ws.Cells.Item(150, 5).Value = 227
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = 18

# Row 167: Guyana
$ws.Cells.Item(167, 2).Value = 290
$ws.Cells.Item(167, 3).Value = 4
$ws.Cells.Item(167, 4).Value = 134
$ws.Cells.Item(167, 5).Value = 140

# Row 175: was Gibraltar, now Lesoto (Lesoto jumps ahead of Gibraltar/Bermudas/Camboya/Brunei)
$ws.Cells.Item(175, 1).Value = "Lesoto"
$ws.Cells.Item(175, 2).Value = 184
$ws.Cells.Item(175, 3).Value = 50
$ws.Cells.Item(175, 4).Value = 26
$ws.Cells.Item(175, 5).Value = 157
$ws.Cells.Item(175, 8).Value = 1

# Row 176: was Bermudas, now Gibraltar (shifted down one rank, same totals as before)
$ws.Cells.Item(176, 1).Value = "Gibraltar"
$ws.Cells.Item(176, 2).Value = 180
$ws.Cells.Item(176, 4).Value = 176
$ws.Cells.Item(176, 5).Value = 4
$ws.Cells.Item(176, 8).Value = 0

# Row 177: was Camboya, now Bermudas (shifted down one rank, same totals as before)
$ws.Cells.Item(177, 1).Value = "Bermudas"
$ws.Cells.Item(177, 2).Value = 149
$ws.Cells.Item(177, 4).Value = 137
$ws.Cells.Item(177, 5).Value = 3
$ws.Cells.Item(177, 8).Value = 9

# Row 178: was Brunei, now Camboya (shifted down one rank, same totals as before)
$ws.Cells.Item(178, 1).Value = "Camboya"
$ws.Cells.Item(178, 4).Value = 131
$ws.Cells.Item(178, 5).Value = 10
$ws.Cells.Item(178, 8).Value = 0

# Row 179: was Lesoto, now Brunei (shifted down one rank, same totals as before)
$ws.Cells.Item(179, 1).Value = "Brunei"
$ws.Cells.Item(179, 2).Value = 141
$ws.Cells.Item(179, 4).Value = 138
$ws.Cells.Item(179, 5).Value = 0
$ws.Cells.Item(179, 8).Value = 3

# Row 189: Antigua y Barbuda
$ws.Cells.Item(189, 2).Value = 74
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 4).Value = 58
$ws.Cells.Item(189, 5).Value = 14
